$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $cellRef, $val) {
    $cell = $ws.Range($cellRef)
    $st = $cell.Style
    $cell.Value = "'" + $val
    $cell.Style = $st
}

Set-TextValue $ws "D2" "34.067.81"
Set-TextValue $ws "E2" "  -0.12%  "
Set-TextValue $ws "D3" "1.788.45"
Set-TextValue $ws "E3" "  -0.06%  "
Set-TextValue $ws "E4" "  +0.06%  "
Set-TextValue $ws "D5" "226.44"
Set-TextValue $ws "E5" "  +1.87%  "
Set-TextValue $ws "E6" "  -1.55%  "
Set-TextValue $ws "E7" "  +0.05%  "
Set-TextValue $ws "D8" "32.25"
Set-TextValue $ws "E8" "  -0.21%  "
Set-TextValue $ws "D9" "0.296"
Set-TextValue $ws "E9" "  +3.65%  "
Set-TextValue $ws "D10" "0.0683"
Set-TextValue $ws "E10" "  -4.43%  "
Set-TextValue $ws "E11" "  +0.90%  "
Set-TextValue $ws "D12" "2.045.90"
Set-TextValue $ws "E12" "  +0.08%  "
Set-TextValue $ws "D13" "11.29"
Set-TextValue $ws "E13" "  +3.38%  "
Set-TextValue $ws "D14" "1.790.78"
Set-TextValue $ws "E14" "  +0.02%  "
Set-TextValue $ws "D15" "34.035.79"
Set-TextValue $ws "E15" "  -0.12%  "
Set-TextValue $ws "D16" "0.620"
Set-TextValue $ws "E16" "  -1.25%  "
Set-TextValue $ws "D17" "4.19"
Set-TextValue $ws "E17" "  +0.22%  "
Set-TextValue $ws "D18" "67.76"
Set-TextValue $ws "E18" "  -0.51%  "
Set-TextValue $ws "D19" "242.68"
Set-TextValue $ws "E19" "  -0.70%  "
Set-TextValue $ws "E20" "  -1.49%  "
Set-TextValue $ws "E21" "  -0.05%  "
Set-TextValue $ws "D22" "10.71"
Set-TextValue $ws "E22" "  -1.07%  "
Set-TextValue $ws "E23" "  -0.44%  "
Set-TextValue $ws "D24" "2.06"
Set-TextValue $ws "E24" "  -2.58%  "
Set-TextValue $ws "D25" "161.78"
Set-TextValue $ws "E25" "  +1.70%  "
Set-TextValue $ws "E26" "  +1.02%  "
Set-TextValue $ws "D27" "16.22"
Set-TextValue $ws "E27" "  -0.95%  "
Set-TextValue $ws "E28" "  +0.07%  "
Set-TextValue $ws "E29" "  +0.19%  "
Set-TextValue $ws "E30" "  +2.41%  "
Set-TextValue $ws "E32" "  -1.15%  "
Set-TextValue $ws "D33" "3.58"
Set-TextValue $ws "E33" "  +2.14%  "
Set-TextValue $ws "E34" "  +1.54%  "
Set-TextValue $ws "D35" "1.398.60"
Set-TextValue $ws "E35" "  -0.10%  "
Set-TextValue $ws "D36" "0.652"
Set-TextValue $ws "E36" "  -0.42%  "
Set-TextValue $ws "B37" "RenderToken"
Set-TextValue $ws "C37" "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue $ws "D37" "2.36"
Set-TextValue $ws "E37" "  +8.66%  "
Set-TextValue $ws "B38" "TrustWalletToken"
Set-TextValue $ws "C38" "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextValue $ws "D38" "1.04"
Set-TextValue $ws "E38" "  -0.92%  "
Set-TextValue $ws "E39" "  +1.41%  "
Set-TextValue $ws "D40" "80.14"
Set-TextValue $ws "E40" "  +0.42%  "
Set-TextValue $ws "E41" "  +0.12%  "
Set-TextValue $ws "D42" "0.920"
Set-TextValue $ws "E42" "  -0.23%  "
Set-TextValue $ws "D43" "13.77"
Set-TextValue $ws "E43" "  +13.96%  "
Set-TextValue $ws "E44" "  -1.00%  "
Set-TextValue $ws "E45" "  +8.71%  "
Set-TextValue $ws "E46" "  +2.91%  "
Set-TextValue $ws "E47" "  +2.95%  "
Set-TextValue $ws "D48" "6.04"
Set-TextValue $ws "E48" "  +2.26%  "
Set-TextValue $ws "D49" "107.70"
Set-TextValue $ws "E49" "  +0.02%  "
Set-TextValue $ws "D50" "1.948.09"
Set-TextValue $ws "E50" "  -0.04%  "
Set-TextValue $ws "E51" "  +0.03%  "
